$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 5 (row 6) - now split into a "slippage" task plus the NRT change note
$ws.Range("C6").Value = "Felhasználókezelés backenden - NRT miatt változtatások"
$ws.Range("D6").Value = "csúszás"

# Week 6 (row 7) - refresh token work, also slipped
$ws.Range("C7").Value = "Felhasználókezelés backenden - refresh token megoldása"
$ws.Range("D7").Value = "csúszás"

# Remaining tasks shift down by one week
$ws.Range("C8").Value = "Felhasználókezelés, kosár kezelése - felület létrehozás - bejelentkezés, profil, rendeléseim képernyő"
$ws.Range("C9").Value = "Shop filters - utánaolvasni, hogyan szokás elkészíteni, backend queryk? Megvalósítása, webes kliensben megvalósítás - shopban filter felület, filter service?"
$ws.Range("C10").Value = "TODO"
$ws.Range("C11").Value = "TODO"
$ws.Range("C12").Value = "Fizetés - utánanézni, milyen lehetőségek vannak, hogyan lehet beépíteni őket, ki lehet-e próbálni őket"
$ws.Range("C13").Value = "Webes fizetés megvalósítása"
$ws.Range("C14").Value = "Kereső optimalizálás - utánanézni, milyen módszerek vannak rá, hogyan érdemes csinálni, implementálni"

# Row 17 "nyár" keeps "Android kliens" instead of what used to be there
$ws.Range("C17").Value = "Android kliens"

# Two new rows appended after the old row 17 content, pushing the remaining
# rows down (row 18 keeps "Android kliens", row 19 is the old "TODO: Android wireframe")
$ws.Range("C18").Value = "Android kliens"
$ws.Range("C19").Value = "TODO: Android wireframe"

# Update the active selection to match the authored state
$ws.Range("C8").Select()
